# Applies the PM_result.xlsx vignette update:
#  - Sheet "PM_result$fields": insert a new "$cens" field row, drop the
#    legacy "$npdata/$itdata" row, and append "$errfile"/"$success" rows.
#  - Sheet "PM_result$methods": append the new method rows ($fit, $auc,
#    $report, $sim, $save, $validate, $step, $opt).
#  - Cosmetic view/selection updates to match the saved workbook state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: PM_result$fields
# ---------------------------------------------------------------------
$fields = $wb.Worksheets.Item("PM_result`$fields")

# Insert a new row 5 (pushes "$pred" and everything below down by one) and
# add the new "$cens" field (B only, no comment in column D).
$fields.Rows.Item(5).Insert()
$fields.Cells.Item(5, 2).Value = '\$cens'

# Remove the old combined "$npdata (class: NPAG, list); $itdata (class:
# IT2B, list)" row. After the insert above it now lives at row 73.
$fields.Rows.Item(73).Delete()

# Append the two new field rows at the bottom of the sheet.
$fields.Cells.Item(75, 1).Value = '\$errfile'
$fields.Cells.Item(75, 2).Value = 'Name of error file if it exists'

$fields.Cells.Item(76, 1).Value = '\$success'
$fields.Cells.Item(76, 2).Value = 'Boolean for successful run'

# Match the saved view/selection state.
$fields.Application.ActiveWindow.ScrollRow = 56
$fields.Range("A77").Select()

# ---------------------------------------------------------------------
# Sheet 2: PM_result$methods
# ---------------------------------------------------------------------
$methods = $wb.Worksheets.Item("PM_result`$methods")

$methods.Cells.Item(5, 1).Value = '\$fit'
$methods.Cells.Item(5, 2).Value = 'Fit data using the model in the PM_result object'

$methods.Cells.Item(6, 1).Value = '\$auc'
$methods.Cells.Item(6, 2).Value = 'Calculate auc by supplying a src, e.g. PM_result$auc("op")'

$methods.Cells.Item(7, 1).Value = '\$report'
$methods.Cells.Item(7, 2).Value = 'Regenerate the report'

$methods.Cells.Item(8, 1).Value = '\$sim'
$methods.Cells.Item(8, 2).Value = 'Simulate using the model in the PM_result object'

$methods.Cells.Item(9, 1).Value = '\$save'
$methods.Cells.Item(9, 2).Value = 'Save the PM_result object'

$methods.Cells.Item(10, 1).Value = '\$validate'
$methods.Cells.Item(10, 2).Value = 'Validate by simuation to create VPC or NPDE as a PM_valid object'

$methods.Cells.Item(11, 1).Value = '\$step'
$methods.Cells.Item(11, 2).Value = 'Stepwise forward/backward linear regression between covariates and model parameter values'

$methods.Cells.Item(12, 1).Value = '\$opt'
$methods.Cells.Item(12, 2).Value = 'Optimal sampling to create a PM_opt object'

$methods.Range("B12").Select()
